$d = $word.ActiveDocument

# --- Paragraph 1: intro paragraph -----------------------------------------
# Insert "via an I2S interface" after "(DAC)" and before the period.
$d.Content.Find.Execute(
    "(DAC). The audio",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "(DAC) via an I2S interface. The audio",
    2)

# "incorporates" -> "is incorporated"
$d.Content.Find.Execute(
    "sample player incorporates with",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "sample player is incorporated with",
    2)

# --- Sound Playback, first paragraph ---------------------------------------
# "playback of multiple independent sounds" -> "multi-channel playback"
$d.Content.Find.Execute(
    "To support playback of multiple independent sounds,",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "To support multi-channel playback,",
    2)

# "latches this accumulated" -> "latches the accumulated"
$d.Content.Find.Execute(
    "latches this accumulated",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "latches the accumulated",
    2)

# --- Sound Playback, last paragraph ----------------------------------------
# "dictates" -> "determines"
$d.Content.Find.Execute(
    "of a sound dictates the address",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "of a sound determines the address",
    2)
